# Collapse the <id>...</id> run-triplet in the tl_p126r <div> block into a
# single run, updating the id value from "p126r_a1" to "p126r_1" in the
# process (matches the "add newly downloaded tc, tcn, tl" commit).
#
# Original runs (3 separate <w:r> elements):
#   [Courier New, color 7f6000, sz 18] "<id>"
#   [default]                          "p126r_a1"
#   [Courier New, color 7f6000, sz 18] "</id>"
#
# Target (single <w:r> element):
#   [Courier New, color 7f6000, sz 18] "<id>p126r_1</id>"

$d = $word.ActiveDocument

# Locate the exact range spanning the old <id>...</id> text (no wildcards,
# match case) so we don't depend on hard-coded character offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("<id>p126r_a1</id>", $true, $false, $false,
                                  $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Re-assigning .Text on a range spanning multiple runs collapses the
    # range into a single run that inherits the formatting of the range's
    # first run (Courier New / 7f6000 / sz 18), exactly as the diff wants.
    $findRange.Text = "<id>p126r_1</id>"
}
